$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.343.93'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '2.355.37'
$ws.Range("E3").Value = '  +5.29%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.44'
$ws.Range("E5").Value = '  +1.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.644'
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.72'
$ws.Range("E7").Value = '  +14.55%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.496'
$ws.Range("E9").Value = '  +13.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0977'
$ws.Range("E10").Value = '  +2.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '27.29'
$ws.Range("E11").Value = '  +1.81%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.106'
$ws.Range("E12").Value = '  +2.47%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.705.19'
$ws.Range("E13").Value = '  +5.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.23'
$ws.Range("E14").Value = '  +5.84%  '
$ws.Range("E15").Value = '  +6.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.868'
$ws.Range("E16").Value = '  +5.74%  '
$ws.Range("D17").Value = '2.367.52'
$ws.Range("E17").Value = '  +5.86%  '
$ws.Range("D18").Value = '43.248.75'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("E19").Value = '  +4.76%  '
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '74.84'
$ws.Range("E20").Value = '  +2.92%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.35'
$ws.Range("E21").Value = '  +5.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '250.35'
$ws.Range("E22").Value = '  +2.33%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  +3.08%  '
$ws.Range("E25").Value = '  +3.20%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.10'
$ws.Range("E26").Value = '  +4.10%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.21'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.42'
$ws.Range("E28").Value = '  +4.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.13'
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.53'
$ws.Range("E30").Value = '  +9.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.132'
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("E32").Value = '  +2.96%  '
$ws.Range("E33").Value = '  +3.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0693'
$ws.Range("E34").Value = '  +3.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.06'
$ws.Range("E35").Value = '  +4.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.74'
$ws.Range("E36").Value = '  +4.32%  '
$ws.Range("E37").Value = '  +4.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.44'
$ws.Range("E38").Value = '  +7.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0254'
$ws.Range("E39").Value = '  +2.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.05'
$ws.Range("E40").Value = '  +12.94%  '
$ws.Range("E41").Value = '  +5.05%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  +1.00%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.38'
$ws.Range("E44").Value = '  +4.03%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.16'
$ws.Range("E45").Value = '  +9.87%  '
$ws.Range("E46").Value = '  +3.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0960'
$ws.Range("E47").Value = '  +2.77%  '
$ws.Range("D48").Value = '1.444.47'
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").Value = '2.580.45'
$ws.Range("E49").Value = '  +5.63%  '
$ws.Range("E50").Value = '  -2.55%  '
$ws.Range("E51").Value = '  +0.75%  '
